# Update the 5x20 "within100" arithmetic-practice table: replace each
# cell's equation text with the corresponding new equation from the
# commit. Cells are visited row-major (Word's natural table order,
# matching $t.Cell(row, col) semantics) and each cell's Range.Text is
# set directly (rather than via Find/Replace on a shared search scope)
# so a replacement can never be mistaken for -- or corrupt -- the text
# of any other, unrelated cell.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @("23+8=", "24+41=", "12-11=", "51-18=", "79-21=", "97-81=", "96-37=", "98-24=", "4+71=", "77+18=", "62-1=", "86-6=", "81+3=", "46-36=", "34-31=", "73-26=", "9+44=", "24+62=", "32+17=", "31+8=", "74-63=", "11+34=", "5+53=", "72+0=", "90-63=", "22+53=", "42+9=", "82-82=", "24+35=", "13+9=", "50-12=", "90+8=", "7+52=", "90-78=", "91-90=", "69-14=", "48+51=", "10+63=", "40-24=", "85+6=", "88-23=", "84-28=", "27+43=", "91-54=", "94+4=", "64-24=", "31+32=", "31+13=", "42-20=", "60-53=", "20+61=", "67-46=", "0+47=", "29+61=", "13+22=", "84-50=", "19+26=", "50-31=", "55+2=", "43-8=", "36+45=", "63-51=", "56-23=", "72-27=", "15+4=", "91+0=", "30-9=", "41-33=", "90-69=", "41+43=", "40+21=", "46-19=", "34+30=", "52+17=", "68+12=", "57+13=", "95-90=", "5+62=", "41+39=", "7+50=", "96-30=", "76-40=", "66-43=", "75+24=", "72-4=", "87-0=", "23+23=", "22+8=", "47+46=", "33+59=", "28+18=", "7+11=", "75-27=", "88-72=", "73-6=", "20-18=", "42+54=", "9+17=", "23+53=", "7+32=")

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Host "Updated $idx cells"
